$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.479.73'
$ws.Range('E2').Value = '  -2.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.620.67'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.14'
$ws.Range('E5').Value = '  -2.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.49'
$ws.Range('E6').Value = '  -4.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.611.46'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -5.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.76'
$ws.Range('E11').Value = '  +15.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.605'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '48.37'
$ws.Range('E13').Value = '  -4.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000284'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.200.78'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '674.40'
$ws.Range('E16').Value = '  -3.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.93'
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.621.40'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.472.54'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.76'
$ws.Range('E21').Value = '  -4.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.45'
$ws.Range('E22').Value = '  -2.69%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.09'
$ws.Range('E24').Value = '  -4.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.62'
$ws.Range('E25').Value = '  -5.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.92'
$ws.Range('E26').Value = '  -2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.78'
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.86'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.59'
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.12'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('E32').Value = '  -5.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.59'
$ws.Range('E33').Value = '  +2.17%  '
$ws.Range('E34').Value = '  -6.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.99'
$ws.Range('E35').Value = '  -5.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '580.29'
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.06'
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('E38').Value = '  -0.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '58.21'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0453'
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.562.95'
$ws.Range('E42').Value = '  -2.44%  '
$ws.Range('E43').Value = '  -3.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.345'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.41'
$ws.Range('E45').Value = '  -4.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0729'
$ws.Range('E46').Value = '  -6.28%  '
$ws.Range('E47').Value = '  -3.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').Value = '  +2.49%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.48'
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('E51').Value = '  -2.58%  '

$wb.Save()
